# Daily attendance processing - 2025-11-29 22:50:04
# Normalize the "Recorded By" (column G) value ordering on the
# "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact-match replacement table for the "Recorded By" column values that
# need their comma-separated entries reordered.
$map = @{
    "System, dnasr281@gmail.com" = "dnasr281@gmail.com, System"
    "admin@admin.com, System" = "System, admin@admin.com"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
    "system, backup@backdoor.com, System" = "backup@backdoor.com, system, System"
}

$lastRow = $ws.UsedRange.Rows.Count
$colG = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cur = $ws.Cells.Item($r, $colG).Value2
    if ($null -ne $cur -and $map.ContainsKey($cur)) {
        $ws.Cells.Item($r, $colG).Value2 = $map[$cur]
    }
}
